$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the shared-string text MODEL_CONDITION -> MODELCONDITION (this cell currently lives at E1,
# before the column shift below moves it to D1)
$ws.Range("E1").Value = "MODELCONDITION"

# Delete column A entirely; this removes the old A-column values (0/8/10 with the bordered style)
# and shifts columns B:F left into A:E, matching the target layout exactly.
$ws.Columns.Item(1).Delete()
